$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.471.16'
$ws.Range("E2").Value = '  -2.51%  '
$ws.Range("D3").Value = '1.746.93'
$ws.Range("E3").Value = '  -2.96%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '323.66'
$ws.Range("E5").Value = '  -0.24%  '
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4440'
$ws.Range("E7").Value = '  +3.42%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3601'
$ws.Range("E8").Value = '  -0.99%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07455'
$ws.Range("E9").Value = '  -1.41%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.01'
$ws.Range("E10").Value = '  -5.98%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.092'
$ws.Range("E11").Value = '  -2.70%  '
$ws.Range("E12").Value = '  +0.13%  '
$ws.Range("E13").Value = '  -5.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.994'
$ws.Range("E14").Value = '  -3.09%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.090'
$ws.Range("E15").Value = '  -3.92%  '
$ws.Range("D16").Value = '1.750.53'
$ws.Range("E16").Value = '  -3.51%  '
$ws.Range("E17").Value = '  -1.11%  '
$ws.Range("E18").Value = '  -1.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06405'
$ws.Range("E19").Value = '  +0.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.73'
$ws.Range("E21").Value = '  -3.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.831'
$ws.Range("E22").Value = '  -2.83%  '
$ws.Range("D23").Value = '27.520.25'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.10'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.106'
$ws.Range("E25").Value = '  -3.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.74'
$ws.Range("E26").Value = '  +1.43%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.30'
$ws.Range("E27").Value = '  -0.63%  '
$ws.Range("D28").Value = '1.950.13'
$ws.Range("E28").Value = '  -3.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.064'
$ws.Range("E29").Value = '  -7.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.05'
$ws.Range("E30").Value = '  -3.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.072'
$ws.Range("E31").Value = '  -9.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.658'
$ws.Range("E32").Value = '  +3.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.08995'
$ws.Range("E33").Value = '  -0.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.459'
$ws.Range("E34").Value = '  -7.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.91'
$ws.Range("E35").Value = '  -7.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02291'
$ws.Range("E36").Value = '  -3.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2077'
$ws.Range("E37").Value = '  -2.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6309'
$ws.Range("E38").Value = '  -3.25%  '
$ws.Range("E39").Value = '  -2.77%  '
$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.896'
$ws.Range("E40").Value = '  -4.71%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.201'
$ws.Range("E41").Value = '  +0.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.388'
$ws.Range("E43").Value = '  -3.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.728'
$ws.Range("E44").Value = '  -2.91%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.20'
$ws.Range("E45").Value = '  -3.49%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5853'
$ws.Range("E46").Value = '  -2.99%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.692'
$ws.Range("E47").Value = '  -0.57%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '120.78'
$ws.Range("E48").Value = '  -3.87%  '
$ws.Range("E49").Value = '  -3.29%  '
$ws.Range("E50").Value = '  -1.19%  '
$ws.Range("E51").Value = '  -1.75%  '
